$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(4, 10, 17, 20, 22, 24, 28, 31, 32, 34, 37, 43, 47, 48, 54, 56, 57, 60, 63, 66, 68, 71, 78, 84, 87, 88, 89, 91, 94, 95, 96, 97, 101, 103, 105, 106, 109, 113, 119, 123, 124, 126, 133, 137, 140, 143, 144, 147, 150, 154, 155, 156, 158, 164, 165, 166, 172, 174, 177, 180, 182, 187, 189, 190)

foreach ($r in $rows) {
    $ws.Range("G$r`:N$r").Value = "N/A"
}
